$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add row 97: LeetCode 206 "Reverse Linked List" ---
$ws.Range("A96:I96").Copy()
$ws.Range("A97:I97").PasteSpecial(-4122)

$ws.Range("A97").Value = 206
$ws.Range("B97").Value = "Reverse Linked List"
$ws.Range("C97").Value = "#linked-list #recursive "
$ws.Range("D97").Value = "easy"
$ws.Range("E97").Value = 2
$ws.Range("F97").Value = 1
$ws.Range("G97").Value = 5
$ws.Range("H97").Value = 45864
$ws.Range("I97").Value = 45864
$ws.Rows.Item(97).RowHeight = 34

# --- Add row 98: LeetCode 3487 "Maximum Unique Subarray Sum After Deletion" ---
$ws.Range("A96:I96").Copy()
$ws.Range("A98:I98").PasteSpecial(-4122)

$ws.Range("A98").Value = 3487
$ws.Range("B98").Value = "Maximum Unique Subarray Sum After Deletion"
$ws.Range("C98").Value = "#array #set"
$ws.Range("D98").Value = "easy"
$ws.Range("E98").Value = 1
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 10
$ws.Range("H98").Value = 45864
$ws.Range("I98").Value = 45864
$ws.Rows.Item(98).RowHeight = 51

# --- Update view: scroll position & selection to match the new bottom of the sheet ---
$excel.ActiveWindow.ScrollRow = 94
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H98:I98").Select()
